# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (want-to-go count) values to the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 255
$ws1.Range("F3").Value = 79
$ws1.Range("F4").Value = 852
$ws1.Range("F5").Value = 527

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 255
$ws4.Range("F3").Value = 79
$ws4.Range("F4").Value = 852
$ws4.Range("F6").Value = 527
